# "update data to 8/23"
# - Record the Gitksan Demo catch (2548 pieces) for week 33 (Aug 13-19) on
#   the "Sockeye FSC and Demo" summary sheet.
# - Append the underlying per-day "demo catches" rows (Gitksan / Beach
#   seine) for Aug 13-19 (2024) that back that weekly total.

$wb = $excel.ActiveWorkbook

# --- Sheet "Sockeye FSC and Demo": fill in week 33 Gitksan Demo catch ---
$wsSummary = $wb.Worksheets.Item("Sockeye FSC and Demo")
$wsSummary.Range("H14").Value = 2548
[void]$wsSummary.Range("H14").Select()

# --- Sheet "demo catches": append new daily rows ---
$wsDemo = $wb.Worksheets.Item("demo catches")

$newRows = @(
    @{ Date = 45517; Nation = "Gitksan"; Gear = "Beach seine"; Pieces = 642 },
    @{ Date = 45518; Nation = "Gitksan"; Gear = "Beach seine"; Pieces = 376 },
    @{ Date = 45519; Nation = "Gitksan"; Gear = "Beach seine"; Pieces = 973 },
    @{ Date = 45520; Nation = "Gitksan"; Gear = "Beach seine"; Pieces = 60 },
    @{ Date = 45521; Nation = "Gitksan"; Gear = "Beach seine"; Pieces = 304 },
    @{ Date = 45522; Nation = "Gitksan"; Gear = "Beach seine"; Pieces = 132 },
    @{ Date = 45523; Nation = "Gitksan"; Gear = "Beach seine"; Pieces = 61 }
)

# Carry the existing row's formatting (date-formatted A column, etc.) down
# onto each freshly appended row instead of assigning a fresh NumberFormat,
# so no new style/numFmt entries get created in styles.xml.
$wsDemo.Range("A23:D23").Copy()

$startRow = 24
$r = $startRow
foreach ($row in $newRows) {
    [void]$wsDemo.Range("A" + $r + ":D" + $r).PasteSpecial(-4122)  # xlPasteFormats
    $wsDemo.Cells.Item($r, 1).Value = $row.Date
    $wsDemo.Cells.Item($r, 2).Value = $row.Nation
    $wsDemo.Cells.Item($r, 3).Value = $row.Gear
    $wsDemo.Cells.Item($r, 4).Value = $row.Pieces
    $r = $r + 1
}
$excel.CutCopyMode = 0

[void]$wsDemo.Range("F33").Select()

$excel.ActiveWorkbook.Save()
